# Powerpoint writer: consolidate text run nodes.
# Merge each word with its trailing space into a single run (the text
# content stays identical; only the <a:r> run-node boundaries change),
# matching the output produced by the updated writer.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title shape: "A" + " " + "slide"  ->  "A " + "slide"
$title = $s.Shapes.Item("Title 1")
$titleRange = $title.TextFrame.TextRange
$titleRange.Characters(1, 2).Text = "A "

# TextBox shape: "Just"+" "+"an"+" "+"image"+" "+"on"+" "+"this"+" "+"side"
# -> "Just "+"an "+"image "+"on "+"this "+"side"
$box = $s.Shapes.Item("TextBox 3")
$boxRange = $box.TextFrame.TextRange
$boxRange.Characters(1, 5).Text = "Just "
$boxRange.Characters(6, 3).Text = "an "
$boxRange.Characters(9, 6).Text = "image "
$boxRange.Characters(15, 3).Text = "on "
$boxRange.Characters(18, 5).Text = "this "
